$wb = $excel.ActiveWorkbook

# ---- Sheet2: remove the helper "scorecard" table (columns K & M) ----
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("K1:K4").ClearContents()
$ws2.Range("M2:M4").ClearContents()

# Apply the bold/blue header style (same as used in Sheet1 header) to B1:I1
$ws2.Range("B1:I1").Font.Bold = $true
$ws2.Range("B1:I1").Font.Size = 14
$ws2.Range("B1:I1").Font.ThemeColor = 4
$ws2.Range("A1:I1").Rows.RowHeight = 18.75

# Autofit columns A:I so column widths match content
$ws2.Range("A1:I20").Columns.AutoFit()

# Update the active cell/selection on Sheet2
$ws2.Range("F10").Select()

# ---- Sheet3: restructure the scorecard summary table ----
$ws3 = $wb.Worksheets.Item("Sheet3")

# Move the values from column C to column B
$ws3.Range("B2").Value = $ws3.Range("C2").Value
$ws3.Range("B3").Value = $ws3.Range("C3").Value
$ws3.Range("B4").Value = $ws3.Range("C4").Value
$ws3.Range("C2:C4").ClearContents()

# Merge and center the header cell across A1:B1
$ws3.Range("A1:B1").Merge()
$ws3.Range("A1:B1").HorizontalAlignment = -4108

# Autofit column A
$ws3.Columns("A:A").AutoFit()

$ws3.Range("E11").Select()

$wb.Save()
